$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M11").Value = 828598.39
$ws.Range("O11").Value = 468005.34

$ws.Range("O12").Value = 56820.53

$ws.Range("O13").Value = 15074.72

$ws.Range("K17").Value = 37564.2

$ws.Range("K25").Value = 12801

$ws.Range("M26").Value = 139910
$ws.Range("N26").Value = 70330
$ws.Range("O26").Value = 69890
